# Add three new login-credential test rows below the existing
# standard_user / secret_sauce pair (rows 2-4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "locked_out_user"
$ws.Range("B2").Value = "secret_sauce"

$ws.Range("A3").Value = "problem_user"
$ws.Range("B3").Value = "secret_sauce"

$ws.Range("A4").Value = "performance_glitch_user"
$ws.Range("B4").Value = "secret_sauce"

# Widen column A to fit the longer usernames, and give column C
# (currently unused) an explicit custom width as well.
# NOTE: the engine stores/round-trips ColumnWidth in whole-pixel
# buckets, so the requested character-width input has to be chosen
# so that, after that pixel rounding, the persisted <col width="..">
# lands on the desired value (30 and ~23.44 respectively).
$ws.Columns.Item(1).ColumnWidth = 29.1
$ws.Columns.Item(3).ColumnWidth = 22.6

# Matches the selection left behind in the saved workbook.
$ws.Range("C14").Select()
